# Commit: "first name, second name -> full name"
#
# The per-participant label cells used to read:
#   ${graph.get("<GROUP>").get(<N>).participant.lastName} ${graph.get("<GROUP>").get(<N>).participant.firstName} ${graph.get("<GROUP>").get(<N>).participant.team? "(" + graph.get("<GROUP>").get(<N>).participant.team + ")" : null}
# and now read:
#   ${graph.get("<GROUP>").get(<N>).participant.fullName} ${graph.get("<GROUP>").get(<N>).participant.team? "(" + graph.get("<GROUP>").get(<N>).participant.team + ")" : null}
# i.e. the separate lastName/firstName tokens were collapsed into a single fullName token
# (everything else about each cell - which bracket slot it labels - stays the same).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value  = '${graph.get("FOUR").get(15).participant.fullName} ${graph.get("FOUR").get(15).participant.team? "(" + graph.get("FOUR").get(15).participant.team + ")" : null}'
$ws.Range("E8").Value  = '${graph.get("THREE").get(7).participant.fullName} ${graph.get("THREE").get(7).participant.team? "(" + graph.get("THREE").get(7).participant.team + ")" : null}'
$ws.Range("B9").Value  = '${graph.get("FOUR").get(14).participant.fullName} ${graph.get("FOUR").get(14).participant.team? "(" + graph.get("FOUR").get(14).participant.team + ")" : null}'
$ws.Range("B11").Value = '${graph.get("FOUR").get(13).participant.fullName} ${graph.get("FOUR").get(13).participant.team? "(" + graph.get("FOUR").get(13).participant.team + ")" : null}'
$ws.Range("E12").Value = '${graph.get("THREE").get(6).participant.fullName} ${graph.get("THREE").get(6).participant.team? "(" + graph.get("THREE").get(6).participant.team + ")" : null}'
$ws.Range("B13").Value = '${graph.get("FOUR").get(12).participant.fullName} ${graph.get("FOUR").get(12).participant.team? "(" + graph.get("FOUR").get(12).participant.team + ")" : null}'
$ws.Range("B15").Value = '${graph.get("FOUR").get(11).participant.fullName} ${graph.get("FOUR").get(11).participant.team? "(" + graph.get("FOUR").get(11).participant.team + ")" : null}'
$ws.Range("E16").Value = '${graph.get("THREE").get(5).participant.fullName} ${graph.get("THREE").get(5).participant.team? "(" + graph.get("THREE").get(5).participant.team + ")" : null}'
$ws.Range("B17").Value = '${graph.get("FOUR").get(10).participant.fullName} ${graph.get("FOUR").get(10).participant.team? "(" + graph.get("FOUR").get(10).participant.team + ")" : null}'
$ws.Range("B19").Value = '${graph.get("FOUR").get(9).participant.fullName} ${graph.get("FOUR").get(9).participant.team? "(" + graph.get("FOUR").get(9).participant.team + ")" : null}'
$ws.Range("E20").Value = '${graph.get("THREE").get(4).participant.fullName} ${graph.get("THREE").get(4).participant.team? "(" + graph.get("THREE").get(4).participant.team + ")" : null}'
$ws.Range("B21").Value = '${graph.get("FOUR").get(8).participant.fullName} ${graph.get("FOUR").get(8).participant.team? "(" + graph.get("FOUR").get(8).participant.team + ")" : null}'
$ws.Range("B23").Value = '${graph.get("FOUR").get(7).participant.fullName} ${graph.get("FOUR").get(7).participant.team? "(" + graph.get("FOUR").get(7).participant.team + ")" : null}'
$ws.Range("E24").Value = '${graph.get("THREE").get(3).participant.fullName} ${graph.get("THREE").get(3).participant.team? "(" + graph.get("THREE").get(3).participant.team + ")" : null}'
$ws.Range("B25").Value = '${graph.get("FOUR").get(6).participant.fullName} ${graph.get("FOUR").get(6).participant.team? "(" + graph.get("FOUR").get(6).participant.team + ")" : null}'
$ws.Range("B27").Value = '${graph.get("FOUR").get(5).participant.fullName} ${graph.get("FOUR").get(5).participant.team? "(" + graph.get("FOUR").get(5).participant.team + ")" : null}'
$ws.Range("E28").Value = '${graph.get("THREE").get(2).participant.fullName} ${graph.get("THREE").get(2).participant.team? "(" + graph.get("THREE").get(2).participant.team + ")" : null}'
$ws.Range("B29").Value = '${graph.get("FOUR").get(4).participant.fullName} ${graph.get("FOUR").get(4).participant.team? "(" + graph.get("FOUR").get(4).participant.team + ")" : null}'
$ws.Range("B31").Value = '${graph.get("FOUR").get(3).participant.fullName} ${graph.get("FOUR").get(3).participant.team? "(" + graph.get("FOUR").get(3).participant.team + ")" : null}'
$ws.Range("E32").Value = '${graph.get("THREE").get(1).participant.fullName} ${graph.get("THREE").get(1).participant.team? "(" + graph.get("THREE").get(1).participant.team + ")" : null}'
$ws.Range("B33").Value = '${graph.get("FOUR").get(2).participant.fullName} ${graph.get("FOUR").get(2).participant.team? "(" + graph.get("FOUR").get(2).participant.team + ")" : null}'
$ws.Range("B35").Value = '${graph.get("FOUR").get(1).participant.fullName} ${graph.get("FOUR").get(1).participant.team? "(" + graph.get("FOUR").get(1).participant.team + ")" : null}'
$ws.Range("E36").Value = '${graph.get("THREE").get(0).participant.fullName} ${graph.get("THREE").get(0).participant.team? "(" + graph.get("THREE").get(0).participant.team + ")" : null}'
$ws.Range("B37").Value = '${graph.get("FOUR").get(0).participant.fullName} ${graph.get("FOUR").get(0).participant.team? "(" + graph.get("FOUR").get(0).participant.team + ")" : null}'

# The author's selection moved from L49 to E39 when saving this workbook.
$ws.Range("E39").Select() | Out-Null
